# Update the Approved/Rejected status for rows 12 and 13 to "Approved"
# and clear the ReasonToReject ("nil") value for those same rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test-Cases")

$ws.Range("I12").Value = "Approved"
$ws.Range("J12").ClearContents()

$ws.Range("I13").Value = "Approved"
$ws.Range("J13").ClearContents()

# Move the active selection to J13 (matches the saved cursor position).
$ws.Range("J13").Select()
